$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cell values (AgTests/AgPosit corrections) ---
$ws.Range("F273").Value = 26824
$ws.Range("G273").Value = 1343
$ws.Range("F279").Value = 42946
$ws.Range("G279").Value = 3100
$ws.Range("F281").Value = 45522
$ws.Range("F282").Value = 47146
$ws.Range("G282").Value = 2842
$ws.Range("F283").Value = 16913
$ws.Range("G283").Value = 1001
$ws.Range("F285").Value = 41846
$ws.Range("G285").Value = 3427
$ws.Range("F286").Value = 55084
$ws.Range("F287").Value = 58727
$ws.Range("F288").Value = 58306
$ws.Range("F289").Value = 62806
$ws.Range("F291").Value = 14929
$ws.Range("G291").Value = 485
$ws.Range("F292").Value = 81799
$ws.Range("G292").Value = 7255
$ws.Range("F293").Value = 81644
$ws.Range("G293").Value = 5743
$ws.Range("F294").Value = 92130
$ws.Range("G294").Value = 4865
$ws.Range("F299").Value = 64544
$ws.Range("G299").Value = 6802
$ws.Range("F300").Value = 71399
$ws.Range("G300").Value = 7033
$ws.Range("F301").Value = 70661
$ws.Range("G301").Value = 5590
$ws.Range("F302").Value = 77043
$ws.Range("G302").Value = 5689
$ws.Range("F306").Value = 70735
$ws.Range("G306").Value = 7136
$ws.Range("F307").Value = 75237
$ws.Range("G307").Value = 6418
$ws.Range("F308").Value = 15806
$ws.Range("G308").Value = 1101
$ws.Range("F309").Value = 74236
$ws.Range("G309").Value = 5260
$ws.Range("F310").Value = 75114
$ws.Range("G310").Value = 3912
$ws.Range("F311").Value = 62463
$ws.Range("G311").Value = 1992
$ws.Range("F312").Value = 26908
$ws.Range("G312").Value = 899
$ws.Range("F313").Value = 71389
$ws.Range("G313").Value = 3208
$ws.Range("F314").Value = 63432
$ws.Range("G314").Value = 3143
$ws.Range("F315").Value = 55814
$ws.Range("G315").Value = 2635
$ws.Range("F316").Value = 49212
$ws.Range("G316").Value = 2225
$ws.Range("F317").Value = 61972
$ws.Range("G317").Value = 2130
$ws.Range("F318").Value = 49108
$ws.Range("G318").Value = 1193
$ws.Range("F319").Value = 41184
$ws.Range("F320").Value = 76017
$ws.Range("G320").Value = 3649
$ws.Range("F321").Value = 90514
$ws.Range("G321").Value = 2798
$ws.Range("F322").Value = 106371
$ws.Range("G322").Value = 2294
$ws.Range("F323").Value = 212186
$ws.Range("G323").Value = 3156
$ws.Range("F324").Value = 232920
$ws.Range("G324").Value = 2655
$ws.Range("F325").Value = 752798
$ws.Range("G325").Value = 6356
$ws.Range("F326").Value = 427738
$ws.Range("G326").Value = 3758
$ws.Range("F327").Value = 238950
$ws.Range("G327").Value = 2885
$ws.Range("F328").Value = 180472
$ws.Range("G328").Value = 2645
$ws.Range("F329").Value = 88706
$ws.Range("G329").Value = 1792
$ws.Range("F330").Value = 70464
$ws.Range("G330").Value = 1972
$ws.Range("F331").Value = 150049
$ws.Range("G331").Value = 2565
$ws.Range("F332").Value = 424048
$ws.Range("G332").Value = 4126
$ws.Range("F333").Value = 258230
$ws.Range("G333").Value = 2684
$ws.Range("F334").Value = 202288
$ws.Range("G334").Value = 3377
$ws.Range("F335").Value = 129006
$ws.Range("G335").Value = 2869
$ws.Range("F336").Value = 99374
$ws.Range("G336").Value = 3138
$ws.Range("F337").Value = 100512
$ws.Range("G337").Value = 2839
$ws.Range("F338").Value = 216272
$ws.Range("G338").Value = 3036
$ws.Range("F339").Value = 625181
$ws.Range("G339").Value = 5347
$ws.Range("F340").Value = 370890
$ws.Range("G340").Value = 3137
$ws.Range("F341").Value = 293132
$ws.Range("G341").Value = 3600
$ws.Range("F342").Value = 185085
$ws.Range("G342").Value = 3126
$ws.Range("F343").Value = 121551
$ws.Range("G343").Value = 2716
$ws.Range("F344").Value = 122458
$ws.Range("G344").Value = 2392
$ws.Range("F345").Value = 250553
$ws.Range("G345").Value = 3087
$ws.Range("F346").Value = 402206
$ws.Range("G346").Value = 3117
$ws.Range("F347").Value = 249092
$ws.Range("G347").Value = 2392
$ws.Range("F348").Value = 2924
$ws.Range("G348").Value = 14

# --- Clear AgTests/AgPosit (F,G) for rows 349-355: data no longer available ---
$ws.Range("F349:G355").ClearContents()

# --- Append new row 356 (st 24. 02. 2021) ---
$ws.Range("A356").Value = 44250
$ws.Range("B356").Value = 298337
$ws.Range("C356").Value = 15160
$ws.Range("D356").Value = 3547
$ws.Range("E356").Value = 6775
